# PlayerPerformance_3795.xlsx update
# - Adds a "Player Info" sheet (before "ODI Batting")
# - Adds an "ODI Batting Extra" sheet (after "ODI Bowling")
# - Renames MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" / "ODI Bowling"
#   and rewrites the full scorecard URLs down to the bare numeric match code
# - Clears the few stray empty INNING_NUMBER cells on "ODI Batting"

$wb = $excel.ActiveWorkbook

$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBowling = $wb.Worksheets.Item("ODI Bowling")

# --- Add "ODI Batting Extra" right after "ODI Bowling" first -----------
# (inserting the "before ODI Batting" sheet first would shift indices and
#  the "after" insert would otherwise land in the wrong place)
$extra = $wb.Worksheets.Add($null, $odiBowling)
$extra.Name = "ODI Batting Extra"

# --- Add "Player Info" right before "ODI Batting" -----------------------
$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

# Worksheet handles captured before an insertion that lands *before* them
# go stale (they track position, not identity) -- re-resolve everything by
# name now that the final sheet order/count is settled.
$playerInfo = $wb.Worksheets.Item("Player Info")
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Item("ODI Batting Extra")

# =========================================================================
# Player Info
# =========================================================================
$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value = $piHeaders[$c - 1]
    $cell.Font.Bold = $true
}

$playerInfo.Cells.Item(2, 1).Value = "'3795"
$playerInfo.Cells.Item(2, 2).Value = "Joshua H Davey"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Medium"

# =========================================================================
# ODI Batting: MATCH_CARD_LINK -> MATCH_CODE (header + full URL -> bare code)
# =========================================================================
$odiBatting.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingDim = $odiBatting.UsedRange
$battingLastRow = $battingDim.Rows.Count
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $odiBatting.Cells.Item($r, 4)
    $url = $cell.Value2
    if ($url -match 'MatchCode=(\d+)') {
        $cell.Value = "'" + $matches[1]
    }
}

# A handful of rows never had an INNING_NUMBER recorded (did-not-bat
# matches) -- their B cell should simply not exist rather than hold an
# empty string.
foreach ($r in @(16, 19, 27)) {
    $odiBatting.Cells.Item($r, 2).ClearContents()
}

# =========================================================================
# ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE (header + full URL -> bare code)
# =========================================================================
$odiBowling.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingDim = $odiBowling.UsedRange
$bowlingLastRow = $bowlingDim.Rows.Count
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $odiBowling.Cells.Item($r, 2)
    $url = $cell.Value2
    if ($url -match 'MatchCode=(\d+)') {
        $cell.Value = "'" + $matches[1]
    }
}

# =========================================================================
# ODI Batting Extra
# =========================================================================
$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $cell = $extra.Cells.Item(1, $c)
    $cell.Value = $extraHeaders[$c - 1]
    $cell.Font.Bold = $true
}

# MATCH_CODE, BATTING_POSITION (number or $null), NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("3637", 8,     "0", "0", "",       "NO"),
    @("3674", 7,     "0", "0", "2.33%",  "NO"),
    @("3675", 7,     "4", "0", "18.10%", "NO"),
    @("3676", 7,     "",  "",  "",       "NO"),
    @("3721", $null, "",  "",  "",       "NO"),
    @("3725", 8,     "5", "0", "24.88%", "YES"),
    @("3733", $null, "",  "",  "",       "NO"),
    @("3753", 8,     "1", "0", "7.75%",  "NO"),
    @("3761", $null, "",  "",  "",       "NO"),
    @("3764", 8,     "0", "0", "0.48%",  "NO"),
    @("3774", $null, "",  "",  "",       "NO"),
    @("3782", 10,    "0", "0", "1.86%",  "NO"),
    @("3787", 8,     "4", "0", "20.00%", "NO"),
    @("3880", 8,     "3", "0", "16.00%", "NO"),
    @("3919", 8,     "",  "",  "",       "NO"),
    @("3979", 4,     "1", "0", "2.44%",  "NO"),
    @("3980", 4,     "1", "0", "3.47%",  "NO"),
    @("4049", $null, "",  "",  "",       "NO"),
    @("4384", 8,     "1", "1", "10.57%", "YES"),
    @("4386", 8,     "0", "0", "",       "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $extra.Cells.Item($r, 1).Value = "'" + $row[0]
    if ($row[1] -ne $null) {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($row[2] -ne "") { $extra.Cells.Item($r, 3).Value = "'" + $row[2] }
    if ($row[3] -ne "") { $extra.Cells.Item($r, 4).Value = "'" + $row[3] }
    if ($row[4] -ne "") { $extra.Cells.Item($r, 5).Value = "'" + $row[4] }
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r++
}
